$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format (one cell at a time -- a comma-joined multi-area
# Range string silently no-ops NumberFormat on some of its areas) on the
# price cells whose new values would otherwise be auto-detected as numbers,
# so they stay text like the rest of the source data.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated coin list values (price + 1h volume change), including the
# FTXToken/Algorand rank swap and SynthetixNetwork -> NEARProtocol replacement.
$ws.Range("D2").Value = "41.471.22"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "2.199.58"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "253.12"
$ws.Range("E5").Value = "  +3.37%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "69.23"
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "0.589"
$ws.Range("E9").Value = "  +7.65%  "
$ws.Range("D10").Value = "38.39"
$ws.Range("E10").Value = "  +7.37%  "
$ws.Range("D11").Value = "0.0950"
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("D12").Value = "58.32"
$ws.Range("E12").Value = "  +0.44%  "
$ws.Range("D13").Value = "7.21"
$ws.Range("E13").Value = "  +7.57%  "
$ws.Range("D14").Value = "0.105"
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("D15").Value = "2.527.41"
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("D16").Value = "0.885"
$ws.Range("E16").Value = "  +5.29%  "
$ws.Range("D17").Value = "14.79"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "2.197.57"
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("D19").Value = "41.384.00"
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("D20").Value = "0.0₃0954"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").Value = "6.27"
$ws.Range("E21").Value = "  +2.21%  "
$ws.Range("D22").Value = "72.02"
$ws.Range("E22").Value = "  -0.69%  "
$ws.Range("D23").Value = "233.20"
$ws.Range("E23").Value = "  -0.76%  "
$ws.Range("D24").Value = "2.10"
$ws.Range("E24").Value = "  +3.20%  "
$ws.Range("D25").Value = "12.05"
$ws.Range("E25").Value = "  +21.87%  "
$ws.Range("D26").Value = "3.90"
$ws.Range("E26").Value = "  +7.72%  "
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("E28").Value = "  +3.47%  "
$ws.Range("E29").Value = "  -2.43%  "
$ws.Range("D30").Value = "170.41"
$ws.Range("E30").Value = "  -0.21%  "
$ws.Range("E31").Value = "  +1.53%  "
$ws.Range("E32").Value = "  +0.97%  "
$ws.Range("D33").Value = "5.60"
$ws.Range("E33").Value = "  +8.00%  "
$ws.Range("D34").Value = "0.122"
$ws.Range("E34").Value = "  -2.59%  "
$ws.Range("D35").Value = "0.0733"
$ws.Range("E35").Value = "  +2.60%  "
$ws.Range("D36").Value = "26.46"
$ws.Range("E36").Value = "  +16.40%  "
$ws.Range("D37").Value = "4.64"
$ws.Range("E37").Value = "  -0.92%  "
$ws.Range("D38").Value = "4.04"
$ws.Range("E38").Value = "  +4.18%  "
$ws.Range("E39").Value = "  +8.74%  "
$ws.Range("D40").Value = "2.24"
$ws.Range("E40").Value = "  -1.98%  "
$ws.Range("D41").Value = "5.80"
$ws.Range("E41").Value = "  -0.58%  "
$ws.Range("D42").Value = "12.13"
$ws.Range("E42").Value = "  +19.57%  "
$ws.Range("D43").Value = "64.12"
$ws.Range("E43").Value = "  -2.77%  "
$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").Value = "5.01"
$ws.Range("E44").Value = "  +1.30%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "0.205"
$ws.Range("E45").Value = "  +7.83%  "
$ws.Range("E46").Value = "  -2.60%  "
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("E49").Value = "  +4.72%  "
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "2.34"
$ws.Range("E51").Value = "  +2.18%  "
